# Add a new "RESOURCELINK" DB parms row to the "Database" sheet, just
# below the header/CFRS row, and make "Database" the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# The existing note on row 17 (ICCS DSN details) needs to end up one row
# further down once the new row is inserted above it. Grab its text now,
# before any row shifting, and remove it so it can be re-added at the
# correct (post-insert) location.
$resourceLinkNote = $ws.Range("A17").Comment
$noteText = $resourceLinkNote.Text()
$resourceLinkNote.Delete()

# Insert a new blank row above the ESBCONF row (current row 3), pushing
# everything below it (including the comment's former home) down by one.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row with the RESOURCELINK DSN entry.
$ws.Range("A3").Value = "RESOURCELINK"
$ws.Range("B3").Value = "CMLIVE"
$ws.Range("C3").Value = "jm08_cml"
$ws.Range("D3").Formula = '=CONCATENATE( "mqsisetdbparms ",ConfigData!$D$4," -n ",A3," -u ",B3," -p ",C3)'
$ws.Range("F3").Formula = '=CONCATENATE( "mqsicvp ",ConfigData!$D$4," -n ",A3)'
$ws.Range("H3").Value = "Attempt to standardise the DSN across environments."

# A3:F3 already inherited the CFRS row's styling from the row insert
# above; H3 is a brand-new column for this row, so give it the same
# look as the equivalent note cell on the GartanRDS row.
$ws.Range("H3").Style = $ws.Range("H20").Style

# Re-attach the ICCS note to its new home, one row below where it used
# to be (row 17 -> row 18) now that the new row has shifted it down.
$ws.Range("A18").AddComment($noteText)

# Make "Database" the active sheet/tab, with F3 selected (matches the
# new row's "mqsicvp" command cell).
$ws.Activate()
$ws.Range("F3").Select()
